$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.761.25"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").Value = "1.648.54"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").Value = "'216.29"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("E6").Value = "  +1.42%  "

$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("E8").Value = "  +1.98%  "

$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("E10").Value = "  +2.44%  "

$ws.Range("D11").Value = "'0.0842"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "1.879.51"
$ws.Range("E12").Value = "  +1.50%  "

$ws.Range("D13").Value = "1.658.99"
$ws.Range("E13").Value = "  +2.59%  "

$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("E15").Value = "  +2.10%  "

$ws.Range("D16").Value = "'65.46"
$ws.Range("E16").Value = "  +0.97%  "

$ws.Range("D17").Value = "26.773.34"
$ws.Range("E17").Value = "  +1.05%  "

$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("D19").Value = "'218.78"
$ws.Range("E19").Value = "  +2.50%  "

$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("E21").Value = "  +1.97%  "

$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("E23").Value = "  +15.78%  "

$ws.Range("E24").Value = "  +2.69%  "

$ws.Range("D25").Value = "'146.62"
$ws.Range("E25").Value = "  -1.28%  "

$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").Value = "'7.12"
$ws.Range("E28").Value = "  +4.34%  "

$ws.Range("D29").Value = "'15.76"
$ws.Range("E29").Value = "  +1.82%  "

$ws.Range("E30").Value = "  +1.80%  "

$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = "  +1.20%  "

$ws.Range("E33").Value = "  +2.52%  "

$ws.Range("D34").Value = "1.280.77"
$ws.Range("E34").Value = "  +4.77%  "

$ws.Range("E35").Value = "  +3.70%  "

$ws.Range("E36").Value = "  +2.25%  "

$ws.Range("E37").Value = "  +3.62%  "

$ws.Range("D38").Value = "'0.536"
$ws.Range("E38").Value = "  +6.33%  "

$ws.Range("D39").Value = "'0.828"
$ws.Range("E39").Value = "  +4.35%  "

$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("E41").Value = "  +2.81%  "

$ws.Range("E42").Value = "  -1.10%  "

$ws.Range("E43").Value = "  +2.47%  "

$ws.Range("D44").Value = "1.789.72"
$ws.Range("E44").Value = "  +1.64%  "

$ws.Range("E45").Value = "  -0.80%  "

$ws.Range("D46").Value = "'59.80"
$ws.Range("E46").Value = "  +9.40%  "

$ws.Range("E47").Value = "  +2.15%  "

$ws.Range("E48").Value = "  +1.27%  "

$ws.Range("E49").Value = "  +4.26%  "

$ws.Range("D50").Value = "'0.0969"
$ws.Range("E50").Value = "  +2.18%  "

$ws.Range("D51").Value = "'0.408"
$ws.Range("E51").Value = "  +0.37%  "
